$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.332.38'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.155.18'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.92%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.96'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.61'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -4.66%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.148.01'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.06%  '

$ws.Range("E9").Value = '  -2.25%  '

$ws.Range("E10").Value = '  -5.45%  '

$ws.Range("E11").Value = '  -3.04%  '

$ws.Range("E12").Value = '  -3.20%  '

$ws.Range("E13").Value = '  -4.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.18'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.676.02'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.88%  '

$ws.Range("E16").Value = '  -2.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.161.70'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.314.74'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.50'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '452.01'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -5.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.89'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.697'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.57'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -5.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.51'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.15'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.56%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.67'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.81'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.69'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.00'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -7.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.08'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.09%  '

$ws.Range("E33").Value = '  -1.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.37'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.75%  '

$ws.Range("E35").Value = '  -6.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.90'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.13'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0694'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0381'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.54%  '

$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '396.28'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -7.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.97'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.111'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.788.45'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -8.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.248'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.12'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.36'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.89'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.14'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.60%  '

$ws.Range("E51").Value = '  -3.70%  '
